$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6
$ws.Range("F4").Value = -1
$ws.Range("F10").Value = -1
$ws.Range("F12").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("F20").Value = 0
$ws.Range("F24").Value = -2
